# Auto-generated Excel COM-interop script
# Applies the crypto price/volume/coin updates described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'53.988.84"
# Row 3
$ws.Cells.Item(3, 4).Value = "'2.252.35"
$ws.Cells.Item(3, 5).Value = '  +2.56%  '
# Row 5
$ws.Cells.Item(5, 4).Value = "'492.33"
$ws.Cells.Item(5, 5).Value = '  +1.21%  '
# Row 6
$ws.Cells.Item(6, 4).Value = "'127.17"
$ws.Cells.Item(6, 5).Value = '  +1.74%  '
# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.19%  '
# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.61%  '
# Row 9
$ws.Cells.Item(9, 4).Value = "'0.0949"
$ws.Cells.Item(9, 5).Value = '  +2.81%  '
# Row 10
$ws.Cells.Item(10, 4).Value = "'0.152"
$ws.Cells.Item(10, 5).Value = '  +2.43%  '
# Row 11
$ws.Cells.Item(11, 4).Value = "'0.324"
$ws.Cells.Item(11, 5).Value = '  +3.12%  '
# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.21%  '
# Row 13
$ws.Cells.Item(13, 4).Value = "'2.650.95"
$ws.Cells.Item(13, 5).Value = '  +2.44%  '
# Row 14
$ws.Cells.Item(14, 4).Value = "'21.70"
$ws.Cells.Item(14, 5).Value = '  +2.90%  '
# Row 15
$ws.Cells.Item(15, 4).Value = "'53.883.51"
$ws.Cells.Item(15, 5).Value = '  +0.75%  '
# Row 16
$ws.Cells.Item(16, 5).Value = '  +0.43%  '
# Row 17
$ws.Cells.Item(17, 4).Value = "'2.252.92"
$ws.Cells.Item(17, 5).Value = '  +2.16%  '
# Row 18
$ws.Cells.Item(18, 4).Value = "'10.00"
$ws.Cells.Item(18, 5).Value = '  +4.81%  '
# Row 19
$ws.Cells.Item(19, 5).Value = '  +2.76%  '
# Row 20
$ws.Cells.Item(20, 4).Value = "'299.75"
$ws.Cells.Item(20, 5).Value = '  +1.90%  '
# Row 21
$ws.Cells.Item(21, 4).Value = "'6.41"
$ws.Cells.Item(21, 5).Value = '  +4.59%  '
# Row 22
$ws.Cells.Item(22, 4).Value = "'1.00"
$ws.Cells.Item(22, 5).Value = '  +0.16%  '
# Row 23
$ws.Cells.Item(23, 5).Value = '  -1.51%  '
# Row 24
$ws.Cells.Item(24, 4).Value = "'61.85"
$ws.Cells.Item(24, 5).Value = '  -1.12%  '
# Row 25
$ws.Cells.Item(25, 5).Value = '  +1.87%  '
# Row 26
$ws.Cells.Item(26, 5).Value = '  +1.56%  '
# Row 27
$ws.Cells.Item(27, 4).Value = "'2.353.66"
$ws.Cells.Item(27, 5).Value = '  +2.58%  '
# Row 28
$ws.Cells.Item(28, 5).Value = '  +1.66%  '
# Row 29
$ws.Cells.Item(29, 4).Value = "'7.04"
$ws.Cells.Item(29, 5).Value = '  +0.50%  '
# Row 30
$ws.Cells.Item(30, 4).Value = "'166.26"
$ws.Cells.Item(30, 5).Value = '  +0.36%  '
# Row 31
$ws.Cells.Item(31, 4).Value = "'1.59"
$ws.Cells.Item(31, 5).Value = '  +0.89%  '
# Row 32
$ws.Cells.Item(32, 4).Value = "'0.0₃0676"
$ws.Cells.Item(32, 5).Value = '  +2.01%  '
# Row 33
$ws.Cells.Item(33, 2).Value = 'Aptos'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(33, 4).Value = "'5.84"
$ws.Cells.Item(33, 5).Value = '  +2.71%  '
# Row 34
$ws.Cells.Item(34, 2).Value = 'USDe'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(34, 4).Value = "'0.999"
$ws.Cells.Item(34, 5).Value = '  +0.10%  '
# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.23%  '
# Row 36
$ws.Cells.Item(36, 4).Value = "'1.07"
$ws.Cells.Item(36, 5).Value = '  +0.20%  '
# Row 37
$ws.Cells.Item(37, 4).Value = "'17.58"
$ws.Cells.Item(37, 5).Value = '  +1.80%  '
# Row 38
$ws.Cells.Item(38, 4).Value = "'0.885"
$ws.Cells.Item(38, 5).Value = '  +7.23%  '
# Row 39
$ws.Cells.Item(39, 4).Value = "'1.18"
$ws.Cells.Item(39, 5).Value = '  +2.45%  '
# Row 40
$ws.Cells.Item(40, 4).Value = "'3.65"
$ws.Cells.Item(40, 5).Value = '  +3.23%  '
# Row 41
$ws.Cells.Item(41, 4).Value = "'35.68"
$ws.Cells.Item(41, 5).Value = '  -0.35%  '
# Row 42
$ws.Cells.Item(42, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(42, 4).Value = "'0.370"
$ws.Cells.Item(42, 5).Value = '  +1.52%  '
# Row 43
$ws.Cells.Item(43, 2).Value = 'Stacks'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(43, 4).Value = "'1.39"
$ws.Cells.Item(43, 5).Value = '  +2.01%  '
# Row 44
$ws.Cells.Item(44, 5).Value = '  +2.17%  '
# Row 45
$ws.Cells.Item(45, 4).Value = "'4.90"
$ws.Cells.Item(45, 5).Value = '  +2.94%  '
# Row 46
$ws.Cells.Item(46, 4).Value = "'124.52"
$ws.Cells.Item(46, 5).Value = '  -0.60%  '
# Row 47
$ws.Cells.Item(47, 4).Value = "'0.0885"
$ws.Cells.Item(47, 5).Value = '  +0.58%  '
# Row 48
$ws.Cells.Item(48, 4).Value = "'0.538"
$ws.Cells.Item(48, 5).Value = '  +1.09%  '
# Row 49
$ws.Cells.Item(49, 2).Value = 'Hedera'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(49, 4).Value = "'0.0481"
$ws.Cells.Item(49, 5).Value = '  +2.19%  '
# Row 50
$ws.Cells.Item(50, 2).Value = 'Bittensor'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(50, 4).Value = "'234.68"
$ws.Cells.Item(50, 5).Value = '  +1.64%  '
# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.49%  '
